$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("J2").Value = "MV-130"
$ws1.Range("J2").Font.Color = 255
